# Weekly update: a new day of Alcachofa price records was added to the
# "Feria Lagunitas de Puerto Montt" sheet. This pushes every existing
# data row (originally rows 5-35) down by one row, and the newest record
# (date 2022-10-25 / serial 44859) is written into the now-empty row 5.
#
# Using Rows(...).Insert() mirrors exactly what Excel does when a user
# inserts a new row above row 5: the used range / <dimension> grows from
# A1:R35 to A1:R36 automatically, and every cell (values + styles, e.g.
# the date-formatted column D) shifts down with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44859
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Española"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "$/caja 30 unidades"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 30
$ws.Range("R5").Value = "Hortaliza"
